$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that need attendance ("p") marked in column H for this week's class
$rows = @(3, 4, 5, 6, 7, 10, 13, 14, 18, 20, 21, 22, 23, 25, 26, 28, 29)

foreach ($r in $rows) {
    $ws.Range("H$r").Value = "p"
}

$ws.Range("F5").Select()
